# Apply PSM caliper tightening update (0.05 -> 0.02) to the summary workbook.

$wb = $excel.ActiveWorkbook

# ---- Sheet 1: 概况 (Overview) ----
$ws1 = $wb.Worksheets.Item("概况")

$ws1.Range("B3").Value = 2830   # 匹配后样本量
$ws1.Range("B4").Value = 1415   # 处理组数 (匹配后)
$ws1.Range("B5").Value = 1415   # 对照组数 (匹配后)
$ws1.Range("B8").Value = 0.02   # 卡尺值

# ---- Sheet 2: 平衡性汇总 (Balance Summary) ----
$ws2 = $wb.Worksheets.Item("平衡性汇总")

# The bias/reduction columns hold numeric-looking values that must stay
# stored as TEXT (as in the source file). Excel's Value auto-detects
# numeric strings and would coerce them to numbers, so we briefly force
# a Text number format while assigning, then clear the format again so
# no extra styling is left behind on the cells.
$textCells = @("C2","D2","C3","D3","C4","D4","C5","D5","C6","D6")
foreach ($cellref in $textCells) {
    $ws2.Range($cellref).NumberFormat = "@"
}

# Row 2: ln_pgdp
$ws2.Range("C2").Value = "-4.78"
$ws2.Range("D2").Value = "111.8"

# Row 3: ln_pop_density
$ws2.Range("C3").Value = "-2.33"
$ws2.Range("D3").Value = "106.0"

# Row 4: industrial_advanced
$ws2.Range("C4").Value = "6.27"
$ws2.Range("D4").Value = "71.3"

# Row 5: ln_fdi
$ws2.Range("C5").Value = "-9.83"
$ws2.Range("D5").Value = "128.3"

# Row 6: ln_road_area
$ws2.Range("C6").Value = "-9.00"
$ws2.Range("D6").Value = "-59.4"

foreach ($cellref in $textCells) {
    $ws2.Range($cellref).ClearFormats()
}

# Row 5: ln_fdi now satisfies the balance criterion
$ws2.Range("E5").Value = "是"
